$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Contoh" -> "Transaction"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contoh")
$ws.Name = "Transaction"

# ---------------------------------------------------------------------------
# 2. Rebuild the sheet contents: drop the old 2-column sample, write the new
#    transaction/payment test-data grid (A1:K7), drop the old hyperlink.
# ---------------------------------------------------------------------------
$ws.Cells.Clear()
$ws.Hyperlinks.Delete()

# --- First-occurrence order matters: it drives the shared-string table
#     order, so cells are written in this specific sequence rather than
#     strict row/column order. ---
$ws.Range("A1").Value = "checkProduct"
$ws.Range("B1").Value = "typePaymentMethod"
$ws.Range("C1").Value = "namePayment"
$ws.Range("D1").Value = "useKinoKoin"
$ws.Range("E1").Value = "typeKinoKoin"
$ws.Range("F1").Value = "kinoKoin"

$ws.Range("B5").Value = "va"
$ws.Range("C2").Value = "mandiri"
$ws.Range("B2").Value = "bt"
$ws.Range("D2").Value = "tidak"
$ws.Range("A2").Value = "Segar Sari c jeruk, 4, ctn"
$ws.Range("D3").Value = "ya"
$ws.Range("E3").Value = "semua"
$ws.Range("E4").Value = "sebagian"

$ws.Range("H1").Value = "potensialSebelum"
$ws.Range("I1").Value = "potensialSesudah"
$ws.Range("J1").Value = "totalKinoKoin"
$ws.Range("K1").Value = "kinoKoinSetelahTerpotong"
$ws.Range("G1").Value = "totalTagihan"

$ws.Range("A5").Value = "cap panda can, 4, ctn"

# --- Remaining cells (values already exist in the shared-string table). ---
$ws.Range("A3").Value = "Segar Sari c jeruk, 4, ctn"
$ws.Range("B3").Value = "bt"
$ws.Range("C3").Value = "mandiri"

$ws.Range("A4").Value = "Segar Sari c jeruk, 4, ctn"
$ws.Range("B4").Value = "bt"
$ws.Range("C4").Value = "mandiri"
$ws.Range("D4").Value = "ya"
$ws.Range("F4").Value = 10000

$ws.Range("C5").Value = "mandiri"
$ws.Range("D5").Value = "tidak"

$ws.Range("A6").Value = "cap panda can, 4, ctn"
$ws.Range("B6").Value = "va"
$ws.Range("C6").Value = "mandiri"
$ws.Range("D6").Value = "ya"
$ws.Range("E6").Value = "semua"

$ws.Range("A7").Value = "cap panda can, 4, ctn"
$ws.Range("B7").Value = "va"
$ws.Range("C7").Value = "mandiri"
$ws.Range("D7").Value = "ya"
$ws.Range("E7").Value = "sebagian"
$ws.Range("F7").Value = 10000

# ---------------------------------------------------------------------------
# 3. Formatting: header row + data grid get a plain (non-hyperlink) Calibri
#    11 font now that the hyperlink styling on column A is gone.
# ---------------------------------------------------------------------------
$ws.Range("A1:F1").Font.Color = $null
$ws.Range("A2:A7").Font.Underline = -4142
$ws.Range("A2:F7").Font.Color = $null

# Column widths roughly matching the new content.
$ws.Columns.Item(1).ColumnWidth = 20.2
$ws.Columns.Item(2).ColumnWidth = 18.3
$ws.Columns.Item(3).ColumnWidth = 12.7
$ws.Columns.Item(4).ColumnWidth = 10.9
$ws.Columns.Item(5).ColumnWidth = 11.7
$ws.Columns.Item(7).ColumnWidth = 12.9
$ws.Columns.Item(8).ColumnWidth = 15.4
$ws.Columns.Item(9).ColumnWidth = 15.2
$ws.Columns.Item(10).ColumnWidth = 12
$ws.Columns.Item(11).ColumnWidth = 22.9

# Page setup: portrait orientation now declared on this sheet.
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Selection / active sheet: Transaction becomes the active tab (previously
#    Sheet1 was last active), with F3 selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("F3").Select()
